# Apply the "Updated symbol list" data refresh to Sheet1.
# Most of the Price (column D) cells are stored as text in the workbook,
# so we force a Text number-format before writing each numeric-looking
# string; otherwise Excel's COM layer would silently coerce the literal
# into a real number and change the cell's stored type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Price-only refreshes (column D) ---
Set-TextValue "D4"  "6.207"
Set-TextValue "D5"  "0.06092"
Set-TextValue "D7"  "6.712"
Set-TextValue "D8"  "1.358"
Set-TextValue "D10" "0.1579"
Set-TextValue "D11" "0.08087"
Set-TextValue "D12" "0.03332"
Set-TextValue "D13" "0.03116"
Set-TextValue "D14" "0.09279"
Set-TextValue "D15" "3.909"
Set-TextValue "D16" "0.001694"
Set-TextValue "D17" "0.04812"
Set-TextValue "D18" "0.0006160"
Set-TextValue "D19" "0.006185"
Set-TextValue "D20" "0.001098"
Set-TextValue "D21" "0.003390"
Set-TextValue "D22" "0.0001499"
Set-TextValue "D23" "3.693"
Set-TextValue "D24" "2.288"
Set-TextValue "D25" "0.3359"
Set-TextValue "D26" "0.1186"
Set-TextValue "D27" "0.0006165"
Set-TextValue "D40" "0.04608"
Set-TextValue "D41" "0.007147"

# --- Row 42/43: the two coins swapped places (CEJI <-> BKEXToken) ---
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1120"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003129"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Remaining price refreshes ---
Set-TextValue "D44" "0.01022"
Set-TextValue "D45" "0.002970"
Set-TextValue "D46" "0.00006024"
Set-TextValue "D48" "0.7500"
Set-TextValue "D49" "0.05810"
